# Insert a new daily price record for "Zanahoria" (Vega Monumental Concepción)
# as row 418, pushing the existing rows 418-506 down to 419-507.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 418:506 down one position by inserting a new row at 418.
$ws.Rows.Item(418).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(418, 1).Value  = 11
$ws.Cells.Item(418, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(418, 3).Value  = "Bíobío"
$ws.Cells.Item(418, 4).Value  = 45275
$ws.Cells.Item(418, 5).Value  = 8
$ws.Cells.Item(418, 6).Value  = 100114013
$ws.Cells.Item(418, 7).Value  = "Zanahoria"
$ws.Cells.Item(418, 8).Value  = "Sin especificar"
$ws.Cells.Item(418, 9).Value  = "Primera"
$ws.Cells.Item(418, 10).Value = 200
$ws.Cells.Item(418, 11).Value = 6000
$ws.Cells.Item(418, 12).Value = 6000
$ws.Cells.Item(418, 13).Value = 6000
$ws.Cells.Item(418, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(418, 15).Value = "Región Metropolitana"
$ws.Cells.Item(418, 16).Value = 300
$ws.Cells.Item(418, 17).Value = 20
$ws.Cells.Item(418, 18).Value = "Hortaliza"
